$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mySheet (4)")
$ws.Activate()

# Update the bounding/reference values that drive the B3:B100 interpolation formulas
$ws.Range("B2").Value = 0.001
$ws.Range("B101").Value = 24

# Move the active selection to B7 (matches the author's saved selection state)
$ws.Range("B7").Select()
